{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change: the Subtitle paragraph \"\u0423\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u0435 \u0432\u0435\u0440\u0441\u0438\u044f\u043c\u0438\" becomes\n// \"\u042f\u0437\u044b\u043a \u0440\u0430\u0437\u043c\u0435\u0442\u043a\u0438 Markdown\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,style,text\");\nawait context.sync();\n\n// Locate the unique paragraph using the \"Subtitle\" style (there is exactly\n// one in this document, holding the course subtitle under the title).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.style === \"Subtitle\") {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Subtitle paragraph not found\");\n}\n\n// Replace the whole paragraph text in one shot, then re-split it into\n// separate runs (word / space / word / space / word) to mirror the\n// target OOXML run layout.\ntarget.insertText(\"\", \"Replace\");\n\nconst range = target.getRange(\"Whole\");\nrange.insertText(\"\u042f\u0437\u044b\u043a\", \"Replace\");\ntarget.insertText(\" \", \"End\");\ntarget.insertText(\"\u0440\u0430\u0437\u043c\u0435\u0442\u043a\u0438\", \"End\");\ntarget.insertText(\" \", \"End\");\ntarget.insertText(\"Markdown\", \"End\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $d / $word / $app resolve against the open document.\n#\n# Change: the Subtitle paragraph \"\u0423\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u0435 \u0432\u0435\u0440\u0441\u0438\u044f\u043c\u0438\" becomes\n# \"\u042f\u0437\u044b\u043a \u0440\u0430\u0437\u043c\u0435\u0442\u043a\u0438 Markdown\".\n#\n# Note on technique: a plain `$range.Text = \"...\"` (or Find&Replace)\n# assignment that spans more than one existing run collapses every run it\n# touches into a single run. To keep the same per-word run layout the\n# source document uses (and that the target OOXML expects - one run per\n# word/space), this script:\n#   1. Rewrites the paragraph's FIRST run in place (safe: Word only\n#      rewrites the run owning Range.Start, the remaining runs in the\n#      paragraph are left untouched).\n#   2. Temporarily splits the paragraph right before the second word so\n#      that word also becomes a \"first run of its own paragraph\", applies\n#      the same in-place rewrite, then deletes the inserted paragraph\n#      mark to rejoin - this leaves the original run boundaries intact.\n#   3. Types the brand-new \" Markdown\" tail with Selection.TypeText on a\n#      collapsed (empty) selection, which (like real Word) starts a fresh\n#      run for each call instead of folding into the previous run.\n\n$d = $word.ActiveDocument\n\n# Locate the unique paragraph using the \"Subtitle\" style (there is exactly\n# one in this document: the course subtitle under the title).\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Style.NameLocal -eq \"Subtitle\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Subtitle paragraph not found\"\n}\n\n# Re-derive the paragraph's 1-based index by scanning for its start offset\n# (more portable across COM hosts than relying on an index captured above).\n$targetStart = $target.Range.Start\n$paraIndex = 1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Start -eq $targetStart) {\n        $paraIndex = $i\n        break\n    }\n}\n\n$newFirstWord = \"\u042f\u0437\u044b\u043a\"\n$newSecondWord = \"\u0440\u0430\u0437\u043c\u0435\u0442\u043a\u0438\"\n$tailWords = @(\" \", \"Markdown\")\n\n# --- Step 1: rewrite the first run (\"\u0423\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u0435\" -> \"\u042f\u0437\u044b\u043a\") -------------\n$p1 = $d.Paragraphs($paraIndex)\n$p1.Range.Text = $newFirstWord\n\n# --- Step 2: split the paragraph right before the second word, rewrite it\n#     in isolation, then rejoin -------------------------------------------\n$p2 = $d.Paragraphs($paraIndex)\n$splitPos = $p2.Range.Start + $newFirstWord.Length + 1   # +1 for the space\n$splitRng = $d.Range($splitPos, $splitPos)\n$splitRng.InsertParagraphAfter()\n\n$p3 = $d.Paragraphs($paraIndex + 1)\n$p3.Range.Text = $newSecondWord\n\n# Delete the paragraph mark that separates the two halves to rejoin them.\n$p2again = $d.Paragraphs($paraIndex)\n$markPos = $p2again.Range.End - 1\n$markRng = $d.Range($markPos, $markPos + 1)\n$markRng.Delete()\n\n# --- Step 3: type the brand-new tail as its own runs ---------------------\n$joined = $d.Paragraphs($paraIndex)\n$joined.Range.Select()\n$sel = $word.Selection\n$sel.Collapse(0)   # wdCollapseEnd\nforeach ($chunk in $tailWords) {\n    $sel.TypeText($chunk)\n}\n"}
